$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = "Wed Nov 01 15:43:30 EDT 2023"
$ws.Range("B3").Value  = "Wed Nov 01 15:43:39 EDT 2023"
$ws.Range("B4").Value  = "Wed Nov 01 15:43:48 EDT 2023"
$ws.Range("B5").Value  = "Wed Nov 01 15:43:57 EDT 2023"
$ws.Range("B6").Value  = "Wed Nov 01 15:44:06 EDT 2023"
$ws.Range("B7").Value  = "Wed Nov 01 15:44:15 EDT 2023"
$ws.Range("B8").Value  = "Wed Nov 01 15:44:24 EDT 2023"
$ws.Range("B9").Value  = "Wed Nov 01 15:44:34 EDT 2023"
$ws.Range("B10").Value = "Wed Nov 01 15:44:45 EDT 2023"
$ws.Range("B11").Value = "Wed Nov 01 15:44:54 EDT 2023"
$ws.Range("B12").Value = "Wed Nov 01 15:45:03 EDT 2023"
$ws.Range("B13").Value = "Wed Nov 01 15:45:12 EDT 2023"
$ws.Range("B14").Value = "Wed Nov 01 15:45:21 EDT 2023"
$ws.Range("B15").Value = "Wed Nov 01 15:45:30 EDT 2023"
$ws.Range("B16").Value = "Wed Nov 01 15:45:40 EDT 2023"
$ws.Range("B17").Value = "Wed Nov 01 15:45:49 EDT 2023"
